# Update the "取得日時" (acquired timestamp) column for all data rows
# on the "ランサーズ" sheet from 2025-11-09 18:23:48 to 2025-11-09 18:29:13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-11-09 18:23:48"
$newValue = "2025-11-09 18:29:13"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2() -eq $oldValue) {
        $cell.Value = $newValue
    }
}
